$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 496, pushing existing rows 496:567 down to 497:568
$ws.Rows("496:496").Insert()

# Populate the newly inserted row 496 with the new data record
$ws.Range("A496").Value = 9
$ws.Range("B496").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C496").Value = "Metropolitana"
$ws.Range("D496").Value = 45034
$ws.Range("D496").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E496").Value = 13
$ws.Range("F496").Value = 100112032
$ws.Range("G496").Value = "Zapallo italiano"
$ws.Range("H496").Value = "Sin especificar"
$ws.Range("I496").Value = "Primera"
$ws.Range("J496").Value = 250
$ws.Range("K496").Value = 5000
$ws.Range("L496").Value = 6000
$ws.Range("M496").Value = 5500
$ws.Range("N496").Value = "`$/caja 50 unidades"
$ws.Range("O496").Value = "Región Metropolitana"
$ws.Range("P496").Value = 110
$ws.Range("Q496").Value = 50
$ws.Range("R496").Value = "Hortaliza"
